$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task: reporte de comprobantes pagados -> en proceso
$ws.Range("A38").Value = "reporte de comprobantes pagados"
$ws.Range("B38").Value = "en proceso"

# New task: revisar cobranzas no imputadas a ningun comprobante -> no comenzado
$ws.Range("A39").Value = "revisar cobranzas no imputadas a ningun comprobante"
$ws.Range("B39").Value = "no comenzado"

# Leave the selection on the last edited cell, matching the author's final state
$ws.Range("C39").Select()
